# Fix elective course scheduling to use common time slots for both sections
$wb = $excel.ActiveWorkbook

# Section_A and Section_B: swap the DS401 (Elective) / Free values between
# Thu 14:00-15:30 (E5) and Wed 15:30-17:00 (D6), since DS401 is being moved
# to the Wed 15:30-17:00 common slot.
$sectionNames = @("Section_A", "Section_B")
foreach ($name in $sectionNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E5").Value = "Free"
    $ws.Range("D6").Value = "DS401 (Elective)"
}

# Elective_Coordination: DS401's row (row 4) moves from Thu 14:00-15:30
# to Wed 15:30-17:00 so it shares a common slot with the other electives.
$coord = $wb.Worksheets.Item("Elective_Coordination")
$coord.Range("B4").Value = "Wed"
$coord.Range("C4").Value = "15:30-17:00"
